$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "ID tổ chức cấp trên(*)  "
$ws.Range("E1").Select()
